$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two blank rows (15 and 16) in the middle of the activity table are
# being removed, which pulls everything below (the last two activity rows
# plus the COUNTA total row) up by two rows.

# Combined height of the soon-to-be-deleted rows, used to shrink the
# decorative "move and size with cells" shapes anchored further down the
# sheet by the same amount (mirrors what Excel does automatically when
# rows are deleted above a floating shape). Compute the absolute target
# heights up front, before the row geometry changes underneath us.
$deletedHeight = $ws.Rows(15).RowHeight + $ws.Rows(16).RowHeight

$targetHeights = @()
for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)
    $targetHeights += ($shp.Height - $deletedHeight)
}

# Delete the two empty rows; everything below shifts up.
$ws.Rows("15:16").Delete()

# Apply the precomputed heights now that the rows are gone.
for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)
    $shp.Height = $targetHeights[$i - 1]
}

# Re-apply the AutoFilter over the new (smaller) table extent.
$ws.AutoFilterMode = $false
$ws.Range("A3:G17").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$wb.Names.Item(1).RefersTo = "=Plan1!`$A`$3:`$G`$17"

# Update the active selection.
$ws.Range("A1:E1").Select()
